$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 15.2
$ws.Range("F2").Value = 14.7
$ws.Range("G2").Value = 14.4
$ws.Range("H2").Value = 14.2

$ws.Range("E3").Value = 17.2
$ws.Range("F3").Value = 18.4
$ws.Range("G3").Value = 19.9
$ws.Range("H3").Value = 21.2

$ws.Range("E4").Value = 20.6
$ws.Range("F4").Value = 20.3
$ws.Range("G4").Value = 20.1
$ws.Range("H4").Value = 19.9

$ws.Range("D5").Value = 5.7
$ws.Range("E5").Value = 5.5
$ws.Range("F5").Value = 5.3
$ws.Range("G5").Value = 5.2
$ws.Range("H5").Value = 5.1

$ws.Range("E6").Value = 37
$ws.Range("F6").Value = 36.2
$ws.Range("G6").Value = 35.7
$ws.Range("H6").Value = 35.6

$ws.Range("F7").Value = 12.1
$ws.Range("G7").Value = 11.7
$ws.Range("H7").Value = 11.4

$ws.Range("E8").Value = 14
$ws.Range("F8").Value = 13.4
$ws.Range("G8").Value = 13
$ws.Range("H8").Value = 12.6

$ws.Range("F9").Value = 32.1
$ws.Range("G9").Value = 31.8
$ws.Range("H9").Value = 31.4

$ws.Range("E10").Value = 47
$ws.Range("F10").Value = 46.1
$ws.Range("G10").Value = 45.7
$ws.Range("H10").Value = 45.1

$ws.Range("E11").Value = 49.8
$ws.Range("F11").Value = 46.1
$ws.Range("G11").Value = 44.5
$ws.Range("H11").Value = 43.3

# Row 12 (Latin America and the Caribbean): remove all data values C:H
$ws.Range("C12:H12").ClearContents()

# Row 13 (Mexico): remove D13, update E:H
$ws.Range("D13").ClearContents()
$ws.Range("E13").Value = 21.7
$ws.Range("F13").Value = 21.8
$ws.Range("G13").Value = 21.6
$ws.Range("H13").Value = 21.2

$ws.Range("E14").Value = 39.2
$ws.Range("F14").Value = 38
$ws.Range("G14").Value = 37.1
$ws.Range("H14").Value = 36.7

$ws.Range("E15").Value = 19.8
$ws.Range("F15").Value = 19.5
$ws.Range("G15").Value = 19.2
$ws.Range("H15").Value = 18.6

$ws.Range("E16").Value = 36.2
$ws.Range("F16").Value = 35.3
$ws.Range("G16").Value = 34.2
$ws.Range("H16").Value = 33.3

$ws.Range("C17").Value = 25.8
$ws.Range("D17").Value = 22.8
$ws.Range("E17").Value = 20.5
$ws.Range("F17").Value = 19.3
$ws.Range("G17").Value = 17.9
$ws.Range("H17").Value = 16.7

$ws.Range("E18").Value = 29.9
$ws.Range("F18").Value = 29.9
$ws.Range("G18").Value = 30
$ws.Range("H18").Value = 30

$ws.Range("E19").Value = 5.9
$ws.Range("G19").Value = 5.5
$ws.Range("H19").Value = 5.4
